$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("Settings")
$wsCredentials = $wb.Worksheets.Item("Credentials")
$wsConstants = $wb.Worksheets.Item("Constants")

# --- Settings sheet: move the "logF_BusinessProcessName" row (was row 4) up to
# row 2, and move the "OrchestratorQueueName" row (was row 2) down to row 6,
# leaving rows 3-5 blank in between. Capture the original values first so the
# two rows can be swapped without clobbering each other.
$origA2 = "OrchestratorQueueName"
$origB2 = "ProcessABCQueue"
$origC2 = "Orchestrator queue Name. The value must match with the queue name defined on Orchestrator."

$origA4 = "logF_BusinessProcessName"
$origB4 = "Framework"
$origC4 = "Logging field which allows grouping of log data of two or more subprocesses under the same business process name"

# Wipe out the two source rows completely (values + formatting) before
# rewriting them so no stray formatted-but-empty cells are left behind.
$wsSettings.Range("A2:C2").Clear()
$wsSettings.Range("A4:C4").Clear()

# New row 2 <- old row 4 content
$wsSettings.Range("A2").Value = $origA4
$wsSettings.Range("B2").Value = $origB4
$wsSettings.Range("C2").Value = $origC4

# New row 6 <- old row 2 content
$wsSettings.Range("A6").Value = $origA2
$wsSettings.Range("B6").Value = $origB2
$wsSettings.Range("C6").Value = $origC2

# --- Constants sheet: new WorkerLimit-style row so it can be edited from
# Config.xlsx directly.
$wsConstants.Range("A13").Value = "BusinessProcess_ParallelExecutions"
$wsConstants.Range("B13").Value = 10

# --- Selections: set the per-sheet remembered selection for each sheet.
$wsSettings.Range("A3").Select() | Out-Null
$wsCredentials.Range("B4").Select() | Out-Null

# Constants becomes the active/selected tab, selected last so it ends up the
# workbook's active sheet on save.
$wsConstants.Select() | Out-Null
$wsConstants.Range("A13").Select() | Out-Null
